$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 stay text cells (they hold text-formatted numbers/percentages)
# even when the new value looks numeric, matching original inlineStr cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.258.25'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '1.790.31'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '337.60'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '0.4534'
$ws.Range('E7').Value = '  +20.38%  '
$ws.Range('D8').Value = '0.3546'
$ws.Range('E8').Value = '  +5.45%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').Value = '1.140'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').Value = '0.07485'
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('D12').Value = '1.004'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '22.34'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').Value = '6.190'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '7.231'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '1.792.83'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('D17').Value = '0.00001085'
$ws.Range('E17').Value = '  +2.91%  '
$ws.Range('D18').Value = '0.06688'
$ws.Range('E18').Value = '  +1.67%  '
$ws.Range('D19').Value = '81.16'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = '1.000'
$ws.Range('D21').Value = '17.13'
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').Value = '6.377'
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('D23').Value = '28.228.48'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').Value = '11.86'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').Value = '2.386'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').Value = '20.45'
$ws.Range('E26').Value = '  +3.38%  '
$ws.Range('D27').Value = '153.74'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').Value = '2.369'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').Value = '1.999.22'
$ws.Range('E29').Value = '  +1.96%  '
$ws.Range('D30').Value = '1.273'
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('D31').Value = '132.24'
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').Value = '4.070'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').Value = '5.855'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').Value = '0.09402'
$ws.Range('E34').Value = '  +7.72%  '
$ws.Range('D35').Value = '0.02371'
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('D36').Value = '12.06'
$ws.Range('E36').Value = '  -2.35%  '
$ws.Range('D37').Value = '0.6643'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.06210'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '0.2153'
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('D40').Value = '5.166'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').Value = '1.484'
$ws.Range('E41').Value = '  +2.69%  '
$ws.Range('D42').Value = '1.209'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('D43').Value = '8.071'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '0.9996'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '14.02'
$ws.Range('E45').Value = '  +2.58%  '
$ws.Range('D46').Value = '3.854'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('D47').Value = '0.6061'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D48').Value = '128.30'
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('D49').Value = '2.020'
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').Value = '0.07083'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').Value = '1.157'
$ws.Range('E51').Value = '  -1.68%  '

# Restore default cell style so style indices match the original workbook
# (only the number format was needed to keep these as text).
$ws.Range("D2:E51").Style = "Normal"

